# Updated 'SelectDataInExcel' method with Apache POI implementation
# Replaces the sample user record in row 2 of the "UserDetails" sheet
# with a newly generated record (Ardath Tromp), mirroring the extra
# test-data rows produced by the Apache POI based generator.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Ardath"
$ws.Range("B2").Value = "Tromp"
$ws.Range("C2").Value = "ArdathTromp01566"
$ws.Range("D2").Value = "34mty8baa4yhl9n"
$ws.Range("E2").Value = "clayton.johnston@gmail.com"
$ws.Range("F2").Value = "1-398-863-1524"
